$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 19, shifting the existing weekly
# records (old rows 19-43) down to rows 21-45, and populate the two
# freshly inserted rows with the new week's data.
$ws.Rows("19:20").Insert()

# Row 19 - new "Especial" quality record for the week of 2022-11-09
$ws.Cells.Item(19, 1).Value2 = 11
$ws.Cells.Item(19, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(19, 3).Value2 = "Bíobío"
$ws.Cells.Item(19, 4).Value2 = 44874
$ws.Cells.Item(19, 5).Value2 = 8
$ws.Cells.Item(19, 6).Value2 = "Fruta"
$ws.Cells.Item(19, 7).Value2 = 100107
$ws.Cells.Item(19, 8).Value2 = "Otros"
$ws.Cells.Item(19, 9).Value2 = 100107002
$ws.Cells.Item(19, 10).Value2 = "Chirimoya"
$ws.Cells.Item(19, 11).Value2 = "Cultivar IV Región"
$ws.Cells.Item(19, 12).Value2 = "Especial"
$ws.Cells.Item(19, 13).Value2 = 50
$ws.Cells.Item(19, 14).Value2 = 26000
$ws.Cells.Item(19, 15).Value2 = 26000
$ws.Cells.Item(19, 16).Value2 = 26000
$ws.Cells.Item(19, 17).Value2 = "`$/bandeja 10 kilos"
$ws.Cells.Item(19, 18).Value2 = "Provincia de Limarí"
$ws.Cells.Item(19, 19).Value2 = 2600
$ws.Cells.Item(19, 20).Value2 = 10

# Row 20 - new "Primera" quality record for the week of 2022-11-09
$ws.Cells.Item(20, 1).Value2 = 11
$ws.Cells.Item(20, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(20, 3).Value2 = "Bíobío"
$ws.Cells.Item(20, 4).Value2 = 44874
$ws.Cells.Item(20, 5).Value2 = 8
$ws.Cells.Item(20, 6).Value2 = "Fruta"
$ws.Cells.Item(20, 7).Value2 = 100107
$ws.Cells.Item(20, 8).Value2 = "Otros"
$ws.Cells.Item(20, 9).Value2 = 100107002
$ws.Cells.Item(20, 10).Value2 = "Chirimoya"
$ws.Cells.Item(20, 11).Value2 = "Cultivar IV Región"
$ws.Cells.Item(20, 12).Value2 = "Primera"
$ws.Cells.Item(20, 13).Value2 = 50
$ws.Cells.Item(20, 14).Value2 = 23000
$ws.Cells.Item(20, 15).Value2 = 23000
$ws.Cells.Item(20, 16).Value2 = 23000
$ws.Cells.Item(20, 17).Value2 = "`$/bandeja 10 kilos"
$ws.Cells.Item(20, 18).Value2 = "Provincia de Limarí"
$ws.Cells.Item(20, 19).Value2 = 2300
$ws.Cells.Item(20, 20).Value2 = 10
